$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Columns("B").Select() | Out-Null

$ws = $wb.Worksheets.Item("DifferentTypes")
$ws.Activate() | Out-Null

$ws.Range("A3").Value = "Country"
$ws.Range("B3").Value = 19
$ws.Range("C3").Value = $true
$ws.Range("D3").Formula = "=B3*B3"

$ws.Range("A4").Value = "Timezone"
$ws.Range("B4").Value = 21
$ws.Range("C4").Value = $false
$ws.Range("D4").Formula = "=B4*B4"

$ws.Range("D5").Select() | Out-Null
